# Updated UI for Estimation for HL and LL
# Populate the Estimation_Sheet with sample CR rows + a Total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - CR1 / Medium complexity example
$ws.Range("A2").Value = "CR1"
$ws.Range("B2").Value = "Medium"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "4"
$ws.Range("D2").Value = 6.74
$ws.Range("E2").Value = 7.08
$ws.Range("F2").Value = "This is test"

# Row 3 - Complex requirement #2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Complex"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = 11.14
$ws.Range("E3").Value = 12.25
$ws.Range("F3").Font.Bold = $false

# Row 4 - Complex requirement #3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Complex"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "5"
$ws.Range("D4").Value = 19.77
$ws.Range("E4").Value = 21.75
$ws.Range("F4").Font.Bold = $false

# Row 5 - Totals
$ws.Range("A5").Font.Bold = $false
$ws.Range("B5").Font.Bold = $false
$ws.Range("C5").Value = "Total"
$ws.Range("D5").Value = 37.65
$ws.Range("E5").Value = 41.08
$ws.Range("F5").Font.Bold = $false
